$d = $word.ActiveDocument

# 1. "репозитория" -> "репозиторий" (word-ending fix: "Ссылка на репозитория на github" ->
#    "Ссылка на репозиторий на github")
$d.Content.Find.Execute("репозитория", $true, $false, $false, $false, $false,
                         $true, 1, $false, "репозиторий", 2)

# 2. Fix "Цель проекта: cоздать ... позволит о весело" -> "Цель проекта: создать ... позволит весело"
#    a) remove the stray "о " before "весело"
$d.Content.Find.Execute("позволит о весело", $true, $false, $false, $false, $false,
                         $true, 1, $false, "позволит весело", 2)

#    b) fix the Latin "c" that starts "создать" to a Cyrillic "с" and drop its
#       English-language mark, scoped tightly to this one occurrence so the
#       many other Latin "c" characters elsewhere in the document are untouched.
$anchor = $d.Content
$anchor.Find.Execute("Цель проекта: c", $true, $false, $false, $false, $false,
                      $true, 1, $false, "", 0)
$scope = $anchor.Duplicate
$scope.Find.Execute("c", $true, $false, $false, $false, $false,
                     $true, 1, $false, "с", 2)
$scope.LanguageID = "ru-RU"

# 3. Remove "оригинальные " before "саундтреки из"
$d.Content.Find.Execute("услышать оригинальные саундтреки из", $true, $false, $false, $false, $false,
                         $true, 1, $false, "услышать саундтреки из", 2)

# 4. Fix typo "поюедителю" -> "победителю"
$d.Content.Find.Execute("поздравления поюедителю.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "поздравления победителю.", 2)
